# Fix invalid characters (escaped apostrophes "\'" and stray non-breaking
# spaces) in several localization strings on the "string" sheet.
#
# The assignments below are intentionally ordered to match the order the
# corrected strings were re-appended to the shared-string table in the
# authoritative edit (Excel appends newly-typed distinct strings to the end
# of the shared string table in the order they are written).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# app_version (fr) - was "Version de l\'application"
$ws.Range("C10").Value = "Version de l'application"

# reset_summary (fr) - was "Supprimer toutes les notes enregistrées et réinitialiser l\'application"
$ws.Range("C85").Value = "Supprimer toutes les notes enregistrées et réinitialiser l'application"

# to_integer (fr) - was "A l\'entier"
$ws.Range("C104").Value = "A l'entier"

# change_class_summary (lu) - had non-breaking spaces
$ws.Range("E19").Value = "Klickt fir är aktuell Klass ze änneren. All gespäichert Notten ginn dobäi geläscht."

# edit_subjects_summary (lu) - had a non-breaking space
$ws.Range("E40").Value = "Klickt fir är Fächer ze beaarbechten"

# G - Human and social sciences (lu) - had a non-breaking space
$ws.Range("E47").Value = "G - Mënschlech an sozial Wëssenschaften"

# ok / "you can always edit ..." (lu) - had non-breaking spaces
$ws.Range("E76").Value = "Dir kënnt är Fächer an aner Optiounen spéider ëmmer nach an den Astellungen beaarbechten"

# reset_summary (lu) - had a non-breaking space and an escaped apostrophe
$ws.Range("E85").Value = "All gespäichert Notten läschen an d'App zeréckzesetzen"

# Restore the selection/viewport that was active when the file was saved.
$ws.Range("F32").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
